$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author filled the FILENAME (A) and MODEL (B) columns down for every
# remaining row of model answers, matching the values already present in
# row 2 ("Rodier Finding" / "gemma3"), and refreshed the BERT score columns.
for ($r = 3; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = "Rodier Finding"
    $ws.Cells.Item($r, 2).Value = "gemma3"
}

# Leave the selection where the author left it when they saved.
$ws.Range("D10").Select() | Out-Null
